$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-empty A5:B5 cells
$ws.Range("A5").Value = 6000
$ws.Range("B5").Value = 6500

# Add new data rows 6-9
$ws.Range("A6").Value = 7500
$ws.Range("B6").Value = 7505

$ws.Range("A7").Value = 7240
$ws.Range("B7").Value = 8520

$ws.Range("A8").Value = 7542
$ws.Range("B8").Value = 6542

$ws.Range("A9").Value = 8542
$ws.Range("B9").Value = 9545

# Row heights for the new rows (8 and 9 slightly shorter per diff: 13.8 vs 14.4)
$ws.Range("A6:B6").RowHeight = 14.4
$ws.Range("A7:B7").RowHeight = 14.4
$ws.Range("A8:B8").RowHeight = 13.8
$ws.Range("A9:B9").RowHeight = 13.8

# Move selection as recorded in the diff
$ws.Range("I21").Select()
